$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.439.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.26%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.156.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.85%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "537.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.23%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.514"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +9.13%  "
$ws.Range("E9").Value = "  +1.13%  "
$ws.Range("E10").Value = "  +3.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.421"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.140"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.698.27"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.83"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000170"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "58.491.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.168.13"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.22"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "377.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +8.46%  "
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.79"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.97%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.515"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.58%  "
$ws.Range("E26").Value = "  +1.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.990"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +14.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0870"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.47%  "
$ws.Range("E30").Value = "  +3.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.93"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.18"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.18"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "160.69"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.24"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.81%  "
$ws.Range("E37").Value = "  +12.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.34"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.656.30"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +10.32%  "
$ws.Range("E40").Value = "  +6.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0681"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.21"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "38.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.708"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.90%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0281"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.198.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.86%  "
$ws.Range("E48").Value = "  +11.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.22"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.979"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.62%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.29"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.38%  "
